$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = '패션명품 쓱세일 NEW SPRING, NEW ITEM'
$ws.Range("B2").Value = '최대15% 할인 쿠폰팩 + 청구할인 혜택'
$ws.Range("C2").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004551'
$ws.Range("D2").Value = '(0306~12) 패션명품 쓱세일 - 혜택'
$ws.Range("E2").Value = '''2023-03-06'
$ws.Range("F2").Value = '''2023-03-12'
$ws.Range("G2").Value = '[''이벤트/쿠폰 > (0306~12) 패션명품 쓱세일 - 혜택'', ''스마일클럽'', ''#혜택 01 패션의류/잡화 전용 할인 혜택 바로보기'', ''#혜택 02 명품 전용 할인 혜택 바로보기'', ''#혜택 04 매일매일 카드청구할인 바로보기'', ''최대 15% 할인! 패션의류/잡화 전용 할인 혜택'', ''패션의류/잡화 대표 브랜드 15% 상품쿠폰'', ''상품쿠폰 15% - 2만원 이상 상품 구매 시 15% 할인, 최대 3만원 2장(패션명품 쓱세일)'', ''쿠폰 받기(기간 내 ID당 1회)'', ''       상품별 15% 할인'', ''       2만원 이상 상품 구매 시 15% 할인, 최대 3만원 (ID당 발급수량 2장)'', ''       쿠폰 적용 방법'', ''쿠폰은 적용상품 수량 1개에 대해 적용됩니다.(동일상품 복수 개 구매 시 1개 수량에만 적용)'', ''* 상품 상세에서 이미 쿠폰을 다운 받으신 경우 본 행사 페이지에서는 쿠폰을 다운 받으실 수 없으며, 본 행사 페이지에서 이미 쿠폰을 다운 받으신 경우에는 상품 상세에서 쿠폰을 다운 받으실 수 없습니다.'', ''패션의류/잡화 대표 브랜드 쿠폰 대상 상품은 검색결과에서도 쉽게 보여요!'', ''백화점 패션의류/잡화 대표 브랜드 15% 상품쿠폰'', ''상품쿠폰 15% - 5만원 이상 상품 구매 시 15% 할인, 최대 3만원 2장'', ''       5만원 이상 상품 구매 시 15% 할인, 최대 3만원 (ID당 발급수량 2장)'', ''최대 15만원 할인! 명품 전용 할인 혜택'', ''상품쿠폰 7% - 10만원 이상 상품 구매 시 7% 할인, 최대 10만원'', ''상품쿠폰 8% - 30만원 이상 상품 구매 시 8% 할인, 최대 15만원'', ''더 많은 상품을 담아 최대 8% 할인 받는 장바구니 쿠폰'', ''* 타 사이트를 통해 방문하신 고객의 경우, 상품 구매 시 쿠폰이 적용되지 않습니다.'', ''장바구니쿠폰 8% - 7만원 이상 구매 시 최대 3만원 할인(ID당 1장 발급)'', ''쿠폰 받기'', ''W컨셉 SPRING SHOP 15% 쿠폰과 함께 W컨셉의 감성을 담은 봄 패션 아이템을 만나보세요!'', ''신세계백화점몰&신세계몰 첫구매 고객이라면  BONUS 쿠폰 하나 더!'', ''* 신세계백화점몰 및 신세계몰 첫구매 고객 대상으로 발급되는 쿠폰입니다.'', ''장바구니쿠폰 12% - 7만원 이상 구매 시 최대 2만원 할인 ID당 1장 발급(첫구매 전용, 패션/명품 전용, 선착순 1만명)'', ''선착순 쿠폰이 마감되었습니다.'', ''카드 혜택도 놓치지 마세요! 매일매일 카드 청구할인'', ''[SSGPAY전용] 삼성카드 8만원 이상 결제 시 8% 청구할인 (일 10만원 한) 자세히 보기'', ''[SSGPAY전용] SSG.COM카드 Edition2 8만원 이상 결제 시 10% 청구할인 (일 10만원 한) 자세히 보기'', ''[SSGPAY전용] 신한카드 8만원 이상 결제 시 8% 청구할인 (일 10만원 ) 자세히 보기'', ''[SSGPAY전용] KB국민카드 8만원 이상 결제 시 8% 청구할인 (일 10만원 ) 자세히 보기'']'
$ws.Range("E2:F2").ClearFormats()

# Row 3
$ws.Range("A3").Value = '패션명품 쓱세일 23 S/S BEST BRAND FOR U'
$ws.Range("B3").Value = '마시모두띠/칼린 외 TOP 브랜드'
$ws.Range("C3").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004552'
$ws.Range("D3").Value = '(0306~12) 패션명품 쓱세일 - 브랜드'
$ws.Range("E3").Value = '''2023-03-06'
$ws.Range("F3").Value = '''2023-03-12'
$ws.Range("G3").Value = '[''이벤트/쿠폰 > (0306~12) 패션명품 쓱세일 - 브랜드'', ''스마일클럽'', ''알도 - 그랜드오픈 기념 단독특가 + 15% 쿠폰 자세히 보기'']'
$ws.Range("E3:F3").ClearFormats()

# Row 4
$ws.Range("A4").Value = '패션명품 쓱세일 매일 오전 10시, 한정수량 타임딜!'
$ws.Range("B4").Value = '쓱세일 봄 스타일링 추천템'
$ws.Range("C4").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004555'
$ws.Range("D4").Value = '(0306~12) 패션명품 쓱세일 - 특가'
$ws.Range("E4").Value = '''2023-03-06'
$ws.Range("F4").Value = '''2023-03-12'
$ws.Range("G4").Value = '[''이벤트/쿠폰 > (0306~12) 패션명품 쓱세일 - 특가'', ''스마일클럽'', ''#혜택 01 매일 오전 10시 초특가 타임딜 바로보기'', ''한정수량 타임딜'', ''타임딜은 오전 10시에 시작됩니다!'', ''오늘의 타임딜이 종료되었습니다.'', ''타임딜은 10시에 오픈됩니다!'', ''금일 타임딜이 종료되었습니다. 감사합니다.'', ''3/6(월) 오후 7시, SSG 단독 발매! 우븐트랙자켓 2종 S머니 추가혜택+이월상품 최대할인 전 구매고객 대상 사은품 까지'', ''3/7(화) 오후 7시, 라이브 동안만 신상품 30%, 이월상품 25~30% 할인, 구매 사은품 증정까지'', ''3/7(화) 오후 8시, 라이브 동안만 신상품 30%, 이월상품 25~30% 할인, 구매 사은품 증정까지'', ''3/7(화) 오후 9시, 라이브 동안만 신상품 30%, 이월상품 25~30% 할인, 구매 사은품 증정까지'', ''3/8(수) 오후 8시, SSG 단독 발매! 우븐트랙자켓 2종 S머니 추가혜택+이월상품 최대할인 전 구매고객 대상 사은품 까지'', ''3/10(금) 오후 6시, 데카브&CC콜렉트 봄 아우터 한정수량 특 무료배송/무료반품 & 추첨 사은품 증정'']'
$ws.Range("E4:F4").ClearFormats()

# Row 5
$ws.Range("A5").Value = '패션 명품 쓱세일  라이브 커머스 혜택'
$ws.Range("B5").Value = '쓱세일 기간 동안 라이브 보고  경품 받기'
$ws.Range("C5").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004593'
$ws.Range("D5").Value = '[패션명품 쓱세일] SSG.LIVE 라이브 구매 이벤트 + 라이브 캘린더 까지'
$ws.Range("E5").Value = '''2023-03-06'
$ws.Range("F5").Value = '''2023-03-12'
$ws.Range("G5").Value = '[''이벤트/쿠폰 > [패션명품 쓱세일] SSG.LIVE 라이브 구매 이벤트 + 라이브 캘린더 까지'', ''스마일클럽'', ''[패션명품 쓱세일] SSG.LIVE 라이브 구매 이벤트 + 라이브 캘린더 까지'', ''[쓱라이브 이벤트 주의사항]'', ''- 패션쓱 SSG.LIVE 에서 패션쓱 SSG.LIVE 제품 구매시 이벤트 응모한 것으로 간주됩니다.'', ''- 이벤트는 모두 주문번호 1개 기준 집계 및 증정됩니다.'', ''- 구매 취소시 이벤트 당첨은 모두 취소됩니다.'', ''- 이벤트에 응모하시는 경우, 아래 개인정보 제공에 동의하신 것으로 간주되며 제공된 정보는 경품 발송후 즉시 파기합니다.'', ''- 사은품 지급 및 이벤트 혜택 제공'', ''- 사은품 지급 및 이벤트 혜택 제공 관련 업무 종료 후 즉시 파기'', ''6. 본 이벤트 경품은 브랜드사 상황에 따라 사전 고지 없이 변경, 대체 될 수 있습니다.'']'
$ws.Range("E5:F5").ClearFormats()

# Row 6
$ws.Range("A6").Value = '어서와 봄, 할인상품 득템해 봄'
$ws.Range("B6").Value = '창고대개방! 봄 맞이 할인 놓치지 마세요.'
$ws.Range("C6").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004717'
$ws.Range("D6").Value = '(3/2~31) 창고대개방! 어서 와 봄, 할인상품 득템해 봄'
$ws.Range("E6").Value = '''2023-03-02'
$ws.Range("F6").Value = '''2023-03-31'
$ws.Range("G6").Value = '[''이벤트/쿠폰 > (3/2~31) 창고대개방! 어서 와 봄, 할인상품 득템해 봄'', ''스마일클럽'', ''(3/2~31) 창고대개방! 어서 와 봄, 할인상품 득템해 봄'']'
$ws.Range("E6:F6").ClearFormats()

# Row 7
$ws.Range("A7").Value = '첫구매는 반값다딜 <풀무원 편>'
$ws.Range("B7").Value = '풀무원 상품 4개 골라 최대 2만원 할인 + 무료배송'
$ws.Range("C7").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004661&domainSiteNo=7018'
$ws.Range("D7").Value = '반값다딜 - 풀무원 편 (3/2~12)'
$ws.Range("E7").Value = '''2023-03-02'
$ws.Range("F7").Value = '''2023-03-12'
$ws.Range("G7").Value = '[''이벤트/쿠폰 > 반값다딜 - 풀무원 편 (3/2~12)'', ''스마일클럽'', ''STEP1 최근 1년간 쓱배송이 처음이라면, 웰컴 쿠폰 받기'', ''Pulmuone 상품쿠폰 50% 4장, 개별 쿠폰 당 최대 5천원 할인'', ''Pulmuone 상품쿠폰 무료배송 1장, 쓱/새벽배송 2만원 이상 구매 시'', ''      최근 1년 간 이마트몰, 트레이더스 쓱배송/점보택배 및 새벽배송 구매 이력이 없는 고객'', ''       쿠폰 발급 기간'', ''       쿠폰 사용 기간'', ''       상품할인쿠폰, 무료배송쿠폰: 2023년 3월 2일(목)~ 2023년 3월 12일 (일)'', ''       쿠폰 사용 조건'', ''-무료배송 쿠폰 : 이마트 쓱배송/새벽배송 상품 2만원 이상 구매시 무료배송'', ''       쿠폰 발급 대상'', ''(50%, 무료배송) [풀무원] 얇은피 꽉찬속 김치만두 400gx2 할인전 : 9,480원 할인후 : 4,740원'', ''(50%, 무료배송) [풀무원] 국산콩 투컵 두부 600g 할인전 : 5,480원 할인후 : 2,740원'', ''(50%, 무료배송) [풀무원] 바릴라 캠핑 기획팩 2,700G 할인전 : 29,800원 할인후 : 24,800원'', ''(50%, 무료배송) [풀무원] 새콤달콤 유부초밥 330g 할인전 : 5,680원 할인후 : 2,840원'', ''(50%, 무료배송) [풀무원] 행복한 동물복지유정란 15구 780g : 8,980원 할인후 : 4,490원'', ''할인 상품 한 눈에 보러가기'']'
$ws.Range("E7:F7").ClearFormats()

# Row 8
$ws.Range("A8").Value = 'Romantic 화이트데이'
$ws.Range("B8").Value = '장바구니 10% 쿠폰+ 달콤한 사랑을 전하세요'
$ws.Range("C8").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004520'
$ws.Range("D8").Value = 'Romantic 화이트데이 - 달콤한 사랑을 마음 가득 담아 전하세요!'
$ws.Range("E8").Value = '''2023-03-06'
$ws.Range("F8").Value = '''2023-03-12'
$ws.Range("G8").Value = '[''이벤트/쿠폰 > Romantic 화이트데이 - 달콤한 사랑을 마음 가득 담아 전하세요!'', ''스마일클럽'', ''화이트데이 10% 쿠폰 받고 쇼핑 GO! 바로보기'', ''#10% 장바구니 쿠폰'', ''사탕만큼 달콤한 10% 할인 쿠폰 바로보기'', ''10% 장바구니 할인 쿠폰'', ''10% 결제쿠폰 (매일 오전 9시 / 선착순 1만명) - 15만원 이상 구매시 최대 1.5만원 할인'', ''      쿠폰 발급 기간'', ''      23/03/06(월) ~ 23/03/12(일), 매일 오전 9시부터 선착순 발급'', ''      쿠폰 사용 기간'', ''      쿠폰 사용 조건'', ''+PLUS 카드 청구할인 혜택도 챙겨가세요!'', ''구호사계절 프리미엄 버킷백, 구호 SSG 단독특가 최대 10%'', ''프라다 外화이트데이 BEST 명품/잡화 프라다/키플링/빅토리아/스톤헨지 外 인기브랜드 ~ 70% 특가!'', ''스와로브스키스와로브스키 주얼리 기프트 추천 크리스탈처럼 반짝이는 주얼리 일부 품목 20%'']'
$ws.Range("E8:F8").ClearFormats()

# Row 9
$ws.Range("A9").Value = 'SSG 스포츠페어'
$ws.Range("B9").Value = '매일 오전 9시 선착순 10%쿠폰'
$ws.Range("C9").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004731'
$ws.Range("D9").Value = 'SSG 스포츠페어'
$ws.Range("E9").Value = '''2023-03-06'
$ws.Range("F9").Value = '''2023-03-12'
$ws.Range("G9").Value = '[''이벤트/쿠폰 > SSG 스포츠페어'', ''스마일클럽'', ''10% 할인쿠폰'', ''할인의 완성'', ''매일 오전 9시 선착순 1만명'', ''스포츠페어 10% 상품쿠폰'', ''상품쿠폰 10% (매일 오전 9시) (선착순 1만명)'', ''쿠폰은 매일 오전 9시 발급됩니다.'', "발급 후 ''My SSG ▶ 쿠폰 ▶ 보유쿠폰'' 에서 확인 가능합니다.", ''23SS UP TO 44%'', ''최대 20% 할인혜택'', ''쿠폰할인 대전'', ''31% OFF'', ''23SS NEW ARRIVAL 자켓/맨투맨/하의 최대 ~10% OFF'', ''23SS NEW 스쿨룩 UPTO~20%★'', ''풋웨어 BEST 에어맥스97/데이브레이크 ~44%'', ''[나이키코리아공식] 남여 인기 신발 모음 ~45%'', ''조던/런닝/트레이닝룩 ~44% SALE!'', ''CLEARENCE 기획전 ~70% OFF + 15% 쿠폰'', ''23 SPRING 점퍼/티셔츠/니트 외 30%OFF 50종 택1'', ''티셔츠/팬츠/스커트 등 신상~이월 최대50%'', ''F/W 최대 40% 할인 특집'', ''[빈폴골프] 남녀 BEST 50 라운딩웨어☆ + SSG 특별 혜택★ + 최대할인쿠폰★'', ''골프존커머스 BEST 클럽/용품 모음전! 23년 신상부터 이월상품 최대할인까지'', ''인도어화 할인전'']'
$ws.Range("E9:F9").ClearFormats()

# Row 10
$ws.Range("A10").Value = '생필품 특가 위크 10% 선착순 할인 쿠폰'
$ws.Range("B10").Value = '화장지부터 세제까지 우리집 필수템 알뜰 쇼핑'
$ws.Range("C10").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004597'
$ws.Range("D10").Value = '생필품 특가 위크'
$ws.Range("E10").Value = '''2023-03-06'
$ws.Range("F10").Value = '''2023-03-12'
$ws.Range("G10").Value = '[''이벤트/쿠폰 > 생필품 특가 위크'', ''스마일클럽'', ''10% 할인 쿠폰'', ''금주의 할인 특가'', ''      쿠폰 적용 대상'', ''매일 9시 3만명 선착순 상품쿠폰 10% SSG.COM'', ''2만원 이상 구매 시 최대 1만원 할인'', ''금일 쿠폰은 모두 소진 되었습니다.'', ''쿠폰은 매일 오전 9시 발급됩니다'', ''       쿠폰 발급 및 사용기간'', "       쿠폰 발급 후 ''My SSG > 쿠폰 > 보유쿠폰'' 에서 확인 가능합니다.", ''       쿠폰은 SSG.COM 회원만 다운 가능하며(간편가입 포함), 비회원 주문시 쿠폰은 적용할 수 없습니다.'', ''♥[모나리자] 생필품 ~20% 할인♥ 휴지/키친타올/미용티슈/물티슈/각티슈/마스크'', ''오슬림밴드/에어팬츠 기저귀 ~2개50%(맘키즈회원전용)'', ''★브랜드할인★ 센트라린 캡슐세탁세제 모음전(아로퓸/카포드/엘에이/향기부스터)'', ''[★40% SALE★]애경 3월 봄이 왔나봄 할인+쿠폰+다다익선+증정'', ''엘경유켈. 세제 4대 브랜드 대통합 할인'', ''[웰라] 프리미엄 살롱 헤어케어 최대 70% off / 샴푸/ 컨디셔너/트리트먼트/헤어팩/헤어오일'', ''★생필품 특가위크★ 헤어/바디 브랜드별 BEST 단독 할인전'', ''★생필품 특가 위크★헤어케어 1+1/2개이상10% 모음'', ''[쓱배송] 3배 빨라진 새치샴푸 런칭 특가전★(함께살수록 할인)'', ''★천만개 판매돌파★ 지아자/크나이프 SSG 단독할인'', ''닥터 브로너스 SSG 할인 +증정혜택'', ''★단독구성+10%쿠폰★ 뉴트로지나/아비노/존슨즈베이비/아비노베이비/클린앤클리어 BEST 특가'', ''★생필품 특가 위크★바디케어 1+1/2개이상10% 모음'', ''세타필 바디로션, 크림  할인 행사'', ''★단독사은품+최대46%할인가★ 리스테린 쿨민트/그린티 외 250ml&750ml 모음전'', ''[서인국마스크][5%다운쿠폰]38년 필터 기업이 만든 귀가 편한 시노텍스 마스크'', ''1+1 / 다다익선 할인 혜택'', ''세제/헤어/바디 ~40% 할인+적립'', ''~30% OFF + 다다익선  + 증정품'', ''인기 제품&기획SET 최대 ~30% OFF'', ''★배송비특가★헤어플러스 브랜드 위크+할인쿠폰+사은품까지!'', ''리뷰이벤트 단독 SSGX살롱클리닉 청담 살롱 새치 ALL 커버 셀프 염색약 전 상품 특가 15% 쿠폰 할인까지'', ''기획전★ 구매왕+다다익선 10%+구매사은품'', ''도브 바디로션 런칭기념행사 [최대~50%할인]'']'
$ws.Range("E10:F10").ClearFormats()

# Row 11
$ws.Range("A11").Value = 'SPRING IS COMING'
$ws.Range("B11").Value = '다가오는 봄, 나에게 주는 선물부터 달콤한 봄 GIFT까지!'
$ws.Range("C11").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004736'
$ws.Range("D11").Value = 'SPRING IS COMING'
$ws.Range("E11").Value = '''2023-03-06'
$ws.Range("F11").Value = '''2023-03-19'
$ws.Range("G11").Value = '[''이벤트/쿠폰 > SPRING IS COMING'', ''스마일클럽'']'
$ws.Range("E11:F11").ClearFormats()

# Row 12
$ws.Range("A12").Value = 'SSG 브랜드 스포트라이트: 코치'
$ws.Range("B12").Value = 'SSG단독 10% 할인 + 구매 사은품 증정'
$ws.Range("C12").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004624'
$ws.Range("D12").Value = 'SSG 브랜드 스포트라이트 - 코치'
$ws.Range("E12").Value = '''2023-03-06'
$ws.Range("F12").Value = '''2023-03-12'
$ws.Range("G12").Value = '[''이벤트/쿠폰 > SSG 브랜드 스포트라이트 - 코치'', ''스마일클럽'', ''상위 20%의 엄선된 고품질 가죽으로 제작된 코치만의 가죽 지갑 최소한의 가공으로 원가죽 무늬와 느낌을 살린 똑같은 가죽이 아닌 하나 뿐인 고유한 상품'', ''      01 SSG 단독 브랜드 스포트라이트 기간 중 10% 할인 혜택'']'
$ws.Range("E12:F12").ClearFormats()

# Row 13
$ws.Range("A13").Value = 'W컨셉 SPRING SHOP#'
$ws.Range("B13").Value = '모한, 르니나 外 봄 아이템 추천 + 15% 쿠폰 혜택'
$ws.Range("C13").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004587'
$ws.Range("D13").Value = '(0306~0312) W컨셉 SPRING SHOP#'
$ws.Range("E13").Value = '''2023-03-06'
$ws.Range("F13").Value = '''2023-03-12'
$ws.Range("G13").Value = '[''이벤트/쿠폰 > (0306~0312) W컨셉 SPRING SHOP#'', ''스마일클럽'', ''SPECIAL OFFER 15% 쿠폰 바로보기'', ''SPECIAL OFFER 15% 상품 쿠폰'', ''15% 상품쿠폰 - 2만원 이상 구매시 최대 2만원 할인'', ''         쿠폰 발급 및 사용 기간'', ''         2만원 이상 상품 구매 시 15% 할인, 최대 2만원 (기간내 ID당 발급수량 2장)'', ''         쿠폰 적용 방법'', ''쿠폰은 적용상품 수량 1개에 대해 적용됩니다. (동일상품 복수 개 구매 시 1개 수량에만 적용)'', ''상품 상세에서 이미 쿠폰을 다운 받으신 경우 행사 페이지에서는 쿠폰을 다운 받으실 수 없으며, 본 행사 페이지에서 이미 쿠폰을 다운 받으신 경우에도 상품 상세에서 쿠폰을 다운받으실 수 없습니다.'']'
$ws.Range("E13:F13").ClearFormats()

# Row 14
$ws.Range("A14").Value = '민지의 장바구니-화이트데이편'
$ws.Range("B14").Value = '디저트 레시피 영상보고 댓글 달면 SSG머니 추첨 증정'
$ws.Range("C14").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004754&domainSiteNo=7018'
$ws.Range("D14").Value = '민지의 장바구니-화이트데이편 (3/6~3/12)'
$ws.Range("E14").Value = '''2023-03-06'
$ws.Range("F14").Value = '''2023-03-12'
$ws.Range("G14").Value = '[''이벤트/쿠폰 > 민지의 장바구니-화이트데이편 (3/6~3/12)'', ''스마일클럽'', ''       My SSG > 이벤트 현황 > 이벤트 참여내역에서 확인해주세요'', ''       ID당 이벤트 기간 내 1회 응모 가능하며, 당첨은 기간 내 ID 당 1회만 가능합니다.'', ''       이벤트 성격과 다른 내용이거나 불쾌감을 줄 수 있는 댓글은 사전 고지 없이 삭제될 수 있습니다.'']'
$ws.Range("E14:F14").ClearFormats()

# Row 15
$ws.Range("A15").Value = '최대 7.2만원 혜택'
$ws.Range("B15").Value = '+ 쓸 때마다 최대 12% 적립'
$ws.Range("C15").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000000858&siteNo=6005&recruitmentPath=L6007001&eventCode=HPG02'
$ws.Range("D15").Value = 'SSG.COM카드 Edition 2 이벤트 안내 페이지'
$ws.Range("E15").Value = '''2022-07-08'
$ws.Range("F15").Value = '''2023-06-28'
$ws.Range("G15").Value = '[''이벤트/쿠폰 > SSG.COM카드 Edition 2 이벤트 안내 페이지'', ''스마일클럽'', ''SSG.COM카드 Edition 2 이벤트 안내 페이지'', ''SSG.COM카드 Edition2는 SSG.COM에서 최대 7.2만원 혜택 + 쓸 때마다 최대 12% 적립'', ''이달의 혜택 01. SSG.COM카드 Edition2 첫결제 시 5.2만원 쿠폰 할인 바로보기'', ''혜택 01. SSG머니 최대 12% 적립 바로보기'', ''혜택 02. 매월 스마일클럽 가입비 3,900원 지원 바로보기'', ''이달의 혜택. SSG.COM에서 SSG.COM카드 Edition2 첫 결제 시 5.2만원 할인'', ''      직전 6개월간 (2022년 9월 1일 부터 2023년 2월 28일) SSG.COM카드 및 SSG.COM카드 Edition2로 결제 이력이 없고 & SSGPAY에 등록된 현대카드가 없는 회원에 한함'', ''첫 결제 쿠폰 이용 방법'', ''② SSGPAY에 등록된 SSG.COM카드 Edition2로'', ''③ SSG.COM에서 기간 내 5만 3천원 이상 첫 결제 시 5.2만원 쿠폰 할인'', ''      본 혜택은 SSG.COM카드 및 SSG.COM카드 Edition2 로 SSG.COM(이마트몰, 신세계몰, 신세계백화점몰 등)에서 직전 6개월간 (2022.08.01~2023.01.31) SSG.COM카드 및 SSG.COM카드 Edition2로 결제 이력이 없고 & SSGPAY에 등록된 현대카드가 없는 회원에 한해 제공됨'', "      본 혜택은 결제 시 [결제방법 > SSGPAY카드] 내 SSG.COM카드 Edition2 선택 시, ''카드할인 최적가'' 추천에 의해 할인 금액이 자동 적용됨. 단, 5만 3천원 이상 첫 결제 시 적용)", ''다운받은 쿠폰은 SSGPAY에 등록된 SSG.COM카드 Edition2로 SSG.COM에서 바로 결제 시 사용 가능합니다.'', ''쿠폰 사용하러 가기'', ''SSGPAY에 등록된 SSG.COM카드 Edition2로 SSG.COM에서 바로 결제 시 사용 가능합니다.'', ''혜택 01. 장 볼 때마다 SSG머니 최대 12% 적립'', ''SSG.COM에서 최대 12%'', ''       스마일클럽으로 5% 적립'', ''       쓱·새벽·트레이더스 구매 시 (구매 전 스마일클럽 적립 아이콘을 꼭 확인해주세요)'', ''       SSG.COM 카드 Edition2로 7% 적립'', ''어디서나 한도없이 0.5%'', ''SSG.COM카드 Edition2로 어디서나 한도없이 0.5% 적립 (SSG.COM 외 모든 가맹점)'', ''SSG Money 최대 12% 적립'', ''스마일클럽 5% 적립 + SSG.COM카드 Edition2 최대 7% 적립'', ''      스마일클럽 5% 적립은 쓱·새벽배송·트레이더스 이용 시에 한함'', ''      SSG.COM카드 Edition2 최대 7% 적립'', ''      SSG.COM에서 결제 시 7%(1만 쓱머니 한도), 그 외 가맹점 0.5% 적립(적립한도 제한 없음)'', ''      무이자 할부 및 현대카드에서 제공하는 다른 할인 서비스 이용 시 적립 제외'', ''스마일클럽 가입비 3,900원 매월 100% 지원'', ''      1. 스마일클럽 자동 가입에'', ''      2. SSGPAY 내 카드 자동 등록'', ''      3. 스마일클럽 정기결제수단 자동 등록 및 월 이용료 3,900원 지원까지! (단, 해당 카드를 월 정기결제 수단에 등록한 경우에 한함)'', ''SSG.COM카드 Edition2를 스마일클럽 월 정기결제 수단에 등록 및 전월 이용금액 30만원 이상 시 혜택 제공'', ''TIP. 스마일클럽 가입 시 SSG.COM 혜택'', ''(SSG 가입 시) 15,000원만 담아도 쓱 · 새벽배송 무료배송'', ''쓱 · 새벽배송 · 트레이더스 SSG머니 최대 5% 적립'', ''장바구니 최대 10% 할인쿠폰'', ''매월 4장씩 최대 12% 할인쿠폰'', ''스마일배송 1만 5천원 이상 무제한 무료배송'', ''스마일배송 상품 스마일캐시 1% 적립'', ''스마일클럽 단독 혜택 스타벅스 상품 전용 딜'', ''스마일클럽 가입비 매월 3,900원 지원'', ''월 1회, 매달 스마일클럽 정기결제일에 혜택 제공'', ''SSG.COM카드 Edition2는 최초 발급 시, 스마일클럽 월 정기결제 수단에 자동 등록 됨'', ''전월 이용금액 30만원 미만 시, SSG.COM카드 Edition2로 스마일클럽 정기 결제 금액이 자동 결제됨'', ''스마일클럽 무료 이용 기간이라면 정기결제 금액 지원 대신 SSG머니 3,900원 제공'', ''      스타벅스 자동 충전, 생활요금(통신요금, 아파트관리비 등) 정기결제 신청 및 이체 시 최대 1만원 청구 할인'', ''      2. 스타벅스 자동 충전 또는 생활요금 정기결제 신청(각 항목당 할인한도 5천원, 최대 1만원 할인)'', ''      이용금액이 혜택금액보다 적을 경우, 이용 금액만큼 할인 적용'', ''      정기결제 신청 후 카드 결제일에 따라 매출 발생 다음 달 또는 다다음 달 청구 할인 혜택 적용'', ''      단, 청구 할인 제공 일정은 당사 또는 신청인 사정에 의해 상이할 수 있음'', ''      3. 쏘카 1만원 할인쿠폰'', ''      쿠폰은 등록일 포함 30일간 이용 가능'', ''실물 SSG.COM카드 Edition2 수령 전 SSGPAY로 결제 시 건당 100만원 이하 결제 가능 *단, 본인 확인(신분증 확인 및 1원 인증) 완료한 경우에 한하며, 건당 100만원 초과 시 실물카드 수령 후 결제 가능'', ''SSG머니 최대 적립 12%에서 5%는 SSG.COM에서 제공하는 멤버십 서비스로 SSG.COM 사정에 따라 변경 가능함'', ''카드 이용대금 연체 시 약정금리 + 연체가산금리 3%의 연체이자율이 적용됩니다. (회원별, 이용 상품별 차등적용 / 법정 최고금리 20% 이내) 단, 연체 발생시점에 약정금리가 없는 경우 아래와 같이 적용'', ''일시불 : 거래 발생시점 기준 최소 기간 (2개월)의 유이자 할부 약정금리 + 연체가산금리 3%'', ''무이자할부 : 거래발생시점 기준 동일한 할부 계약 기간의 유이자할부 약정금리 + 연체가산금리 3%'']'
$ws.Range("E15:F15").ClearFormats()

# Row 16
$ws.Range("A16").Value = '더 강력해진 SSG.COM 삼성카드'
$ws.Range("B16").Value = 'SSG MONEY 최대 15% 적립 + 스마일클럽 월이용료 할인'
$ws.Range("C16").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000002385&recruitmentPath=SSG'
$ws.Range("D16").Value = 'SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지'
$ws.Range("E16").Value = '''2022-10-26'
$ws.Range("F16").Value = '''2025-10-25'
$ws.Range("G16").Value = '[''이벤트/쿠폰 > SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지'', ''스마일클럽'', ''SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지'', ''SSG머니 최대 15% 적립'', ''스마일클럽 월 이용료 3,900원 결제일 할인'', ''최대 15% SSG머니 적립 중 5%는 SSG.COM에서 제공하는 서비스로 자세한 내용은 SSG.COM 확인바람.'', ''01. SSG MONEY 최대 15% 적립!'', ''       이용실적 관계없이 적립 : 1 % + 이용실적 충족 시 적립 4% = 적립 가능한 최대 SSG머니 5%'', ''       카드혜택 + 스마일클럽 혜택'', ''       이용실적 관계없이 적립 1% + 이용실적 충족 시 적립 9% + 쓱배송/쓱배송 traders/새벽배송 상품 구매 시 5% = 적립 가능한 최대 SSG머니 15%'', ''최대 15% SSG MONEY 적립'', ''SSG.COM 삼성카드 최대 10% 적립+스마일클럽 5% 적립'', ''SSG.COM 삼성카드 최대 10% 적립(1% 적립+9% 추가 적립)'', ''1% 적립(전월 실적 조건 및 적립한도 없음)'', ''9% 추가 적립'', ''삼성카드 할인이 적용된 일시불 및 할부 이용금액은 제외됩니다.'', ''스마일클럽 5% 적립'', ''쓱배송/쓱배송 traders/새벽배송 상품 결제건에 한해 혜택을 받을 수 있습니다.'', ''02. 스마일클럽 월이용료 매월 3,900원 결제일 할인'', ''15,000원만 담아도 쓱 새벽배송 무료배송(SSG가입시)'', ''1쓱 새벽배송 트레이더 - SSG머니 최대 5% 적립'', ''장바구니 최대 10% 할인쿠폰'', ''스마일클럽 단독 혜택 - 스타벅스 상품 전용 딜'', ''매월 4장씩 최대 12% 할인쿠폰'', ''스마일배송 1만 5천원 이상 무제한 무료배송'', ''스마일배송 상품 스마일캐시 1% 적립'', ''카드 할인 혜택 자세히보기 (레이어팝업 열기)'', ''2022.10.26부터, SSG.COM 삼성카드 발급 시 스마일클럽에 자동 가입됩니다.'', ''G마켓 또는 옥션을 통해 스마일클럽에 가입한 경우 혜택을 받을 수 없습니다.'', ''SSG.COM 삼성카드로 스마일클럽 월 이용료(3,900원) 정기결제 시 혜택이 제공됩니다.(월 1회)'', ''SSG.COM을 통한 스마일클럽 가입건에 한해 혜택이 제공됩니다.'', ''결제금액이 할인금액보다 적을 경우, 결제금액만큼 할인이 적용됩니다.(결제금액이 없는 경우 할인 대상에서 제외)'', ''01. 5만 5천원 이상 결제 시 사용 가능한 5만원 할인쿠폰 제공'', ''50,000원 할인쿠폰 - SSGPAY 바로결제 이용 시 사용 가능'', ''쿠폰 발급기간 : 2023.03.01 ~ 2023.03.31'', ''쿠폰 사용기간 : 2023.03.01 ~ 2023.03.31'', ''SSGPAY 바로결제 이용 이력이 없는 회원'', ''SSGPAY 바로결제에 등록된 SSG.COM 삼성카드로 결제해야 쿠폰을 사용할 수 있습니다.'', ''배송비 등을 제외한 최종 결제금액이 55,000원 이상이여야 쿠폰을 사용할 수 있습니다.'', ''SSGPAY 바로결제에 등록된 SSG.COM 삼성카드로 결제하셔야 쿠폰을 사용할 수 있습니다.'', ''쿠폰은 통합 회원 본인 명의의 SSG.COM 삼성카드로 결제 시 사용 가능하며, 다른 부정적인 방법으로 사용한 경우에는 주문이 취소될 수 있습니다.'', ''쿠폰은 결제 화면에서 자동으로 적용됩니다.'', ''결제 화면에서 쿠폰 변경을 원할 경우 ‘쿠폰선택’을 눌러주세요.'', ''02. 삼성카드 결제일할인 행사일에는 최대 2% 즉시할인'', ''최대 2% 즉시할인 : [카드혜택가 > SSG.COM 삼성카드 행사 더보기] 내 즉시할인가가 표시된 상품 구매 시 혜택을 받을 수 있습니다. (상품별 할인한도 30,000원)'', ''삼성카드 결제일할인 행사일에만 적용됩니다.'', ''삼성카드 할인이 적용된 결제건은 SSG머니 적립 대상에서 제외됩니다.'', ''03. 스마일클럽 월 이용료 결제 카드로 SSG.COM 삼성카드 등록 시 SSG머니 3,900원 적립'', ''행사기간 동안 SSG.COM 삼성카드를 통해 스마일클럽 자동 가입 시 SSG머니 3,900원 즉시 적립'', ''적립시점까지 스마일클럽 월 이용료 정기결제 수단에 SSG.COM 삼성카드를 등록해야 혜택을 받을 수 있습니다.'', ''본 상품 발급 및 SSGPAY앱 설치시 바로결제에 자동 등록 됩니다.'', ''SSG MONEY는 매월 1일~말일까지 매출전표가 접수된 금액에 대해 다음달 25일 SSG.COM 계정으로 적립됩니다.'', ''적립된 SSG MONEY는 SSGPAY 회원가입 후 조회 및 사용 가능합니다. 단, SSG.COM 에서는 SSG.COM만 가입해도 사용 가능합니다.'', ''연체이자율 : 회원별/이용상품별 정상이자율+3.0%p(최고 연 20.0%)'', ''이미 SSG.COM 삼성카드를 가지고 계시네요!스마일클럽 가입하고 모든 혜택 누리세요'', ''이미 SSG.COM 삼성카드를 가지고 계시네요!G마켓 또는 옥션을 통해 스마일클럽에 가입한 경우, 월 이용료 결제일할인 혜택을 받을 수 없습니다.'', ''쿠폰 다운되었습니다.다운받은 쿠폰은 SSG.COM 삼성카드 발급 후 바로 사용 가능합니다.'', ''이벤트 쿠폰이 이미 발급되었습니다.다운받은 쿠폰은 SSG.COM 삼성카드 발급 후 바로 사용 가능합니다.'']'
$ws.Range("E16:F16").ClearFormats()

# Row 17
$ws.Range("A17").Value = '프레드페리 3/6(월) 7PM'
$ws.Range("B17").Value = '우븐트랙자켓 단독 발매+S머니 할인! 봄신상 & 이월 최대 56%'
$ws.Range("C17").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004605&siteNo=6005'
$ws.Range("D17").Value = '프레드페리@SSG.LIVE 3/6(월) 19:00PM'
$ws.Range("E17").Value = '''2023-02-27'
$ws.Range("F17").Value = '''2023-03-09'
$ws.Range("G17").Value = '[''이벤트/쿠폰 > 프레드페리@SSG.LIVE 3/6(월) 19:00PM'', ''스마일클럽'', ''- 사은품 지급 및 이벤트 혜택 제공'', ''- 사은품 지급 및 이벤트 혜택 제공 관련 업무 종료 후 즉시 파기'']'
$ws.Range("E17:F17").ClearFormats()

# Row 18
$ws.Range("A18").Value = '[월간호캉쓱] 파라다이스 시티 3/6(월) 8PM'
$ws.Range("B18").Value = '1박 + 6만 리조트머니 27만원대~, 2박 + 조식 + 3만 리조트 머니 54만원대~'
$ws.Range("C18").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004712&siteNo=6005'
$ws.Range("D18").Value = '파라다이스 시티 @SSG.LIVE 3/6(월) 20:00'
$ws.Range("E18").Value = '''2023-02-28'
$ws.Range("F18").Value = '''2023-03-06'
$ws.Range("G18").Value = '[''이벤트/쿠폰 > 파라다이스 시티 @SSG.LIVE 3/6(월) 20:00'', ''스마일클럽'']'
$ws.Range("E18:F18").ClearFormats()

# Row 19
$ws.Range("A19").Value = '킨록by킨록앤더슨 3/7(화) 19PM'
$ws.Range("B19").Value = '신상 30% 선할인&이월 ~30% 할인 라이브 한정 특가 진행'
$ws.Range("C19").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004588&siteNo=6005'
$ws.Range("D19").Value = '킨록by 라이브 동안 신상품 30% 할인 @SSG.LIVE 3//7(화) 19PM'
$ws.Range("E19").Value = '''2023-03-01'
$ws.Range("F19").Value = '''2023-03-07'
$ws.Range("G19").Value = '[''이벤트/쿠폰 > 킨록by 라이브 동안 신상품 30% 할인 @SSG.LIVE 3//7(화) 19PM'', ''스마일클럽'', ''킨록by 라이브 동안 신상품 30% 할인 @SSG.LIVE 3//7(화) 19PM'']'
$ws.Range("E19:F19").ClearFormats()

# Row 20
$ws.Range("A20").Value = '[패션명품 쓱세일] 구해줘명품즈 시즌2 3/7(화) 8PM'
$ws.Range("B20").Value = '골든구스/보테가베네타/마르지엘라 외 SS명품 한정수량 역대급 파격 특가'
$ws.Range("C20").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004703&domainSiteNo=6005'
$ws.Range("D20").Value = '구해줘명품즈2 @SSG.LIVE 3/7(화) 20:00PM'
$ws.Range("E20").Value = '''2023-02-27'
$ws.Range("F20").Value = '''2023-03-07'
$ws.Range("G20").Value = '[''이벤트/쿠폰 > 구해줘명품즈2 @SSG.LIVE 3/7(화) 20:00PM'', ''스마일클럽'', ''                    SSG.LIVE 사은품 지급 및 이벤트 혜택 당첨 주의사항'', '' - 사은품 지급 및 이벤트 혜택 제공'', '' - 사은품 지급 및 이벤트 혜택 제공 관련 업무 종료 후 즉시 파기'']'
$ws.Range("E20:F20").ClearFormats()

# Row 21
$ws.Range("A21").Value = '온앤온 3/7(화) 9PM'
$ws.Range("B21").Value = '봄 신상 & 이월 최대 67%! 싱글 투버튼 자켓 선착순 9만9천원 핫딜'
$ws.Range("C21").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004705&siteNo=6005'
$ws.Range("D21").Value = '온앤온 @SSG.LIVE 3/7(화) 9PM'
$ws.Range("E21").Value = '''2023-02-27'
$ws.Range("F21").Value = '''2023-03-08'
$ws.Range("G21").Value = '[''이벤트/쿠폰 > 온앤온 @SSG.LIVE 3/7(화) 9PM'', ''스마일클럽'', ''[SSG.LIVE 사은품 지급 및 이벤트 혜택 당첨 주의사항]'', ''- 사은품 지급 및 이벤트 혜택 제공'', ''- 사은품 지급 및 이벤트 혜택 제공 관련 업무 종료 후 즉시 파기'']'
$ws.Range("E21:F21").ClearFormats()

